$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Step1 - Input Data")
$ws2 = $wb.Worksheets.Item("Step2 - Projection")

# 'Step1 - Input Data'!B4 : Personal Required Rate of Return, stored as literal
# text "8.72%" -> "8.66%" even though the cell carries a "0%" number format.
# Flip to a Text format before the assignment so the percent-looking string
# is kept verbatim instead of being auto-parsed into a numeric percentage,
# then restore the original "0%" number format (re-formatting an existing
# text cell does not re-parse its content).
$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "8.66%"
$ws1.Range("B4").NumberFormat = "0%"

# 'Step2 - Projection'!C5 : Growth Rate, same literal-text-in-a-percent-cell
# situation, "17.15%" -> "17.0%".
$ws2.Range("C5").NumberFormat = "@"
$ws2.Range("C5").Value = "17.0%"
$ws2.Range("C5").NumberFormat = "0%"

# 'Step2 - Projection'!E11 : hard-coded Revenue projection input.
$ws2.Range("E11").Value = 85710

$wb.Application.CalculateFullRebuild()
